$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------------
# 1. Remove the "literal mapping" sheet entirely (folded into the generic
#    `mapping` model — the `literal` slot is handled elsewhere now).
# ---------------------------------------------------------------------------
$litMapping = $wb.Worksheets.Item("literal mapping")
$litMapping.Delete()

# ---------------------------------------------------------------------------
# 2. `mapping` sheet: semantic_similarity_* columns renamed to similarity_*.
# ---------------------------------------------------------------------------
$mapping = $wb.Worksheets.Item("mapping")
$mapping.Range("AM1").Value = "similarity_score"
$mapping.Range("AN1").Value = "similarity_measure"

# ---------------------------------------------------------------------------
# 3. `mapping set` sheet: new `curie_map` slot becomes the first column,
#    pushing every other column one to the right.
# ---------------------------------------------------------------------------
$mappingSet = $wb.Worksheets.Item("mapping set")
$mappingSet.Columns("A:A").Insert()
$mappingSet.Range("A1").Value = "curie_map"

# ---------------------------------------------------------------------------
# 4. New `prefix` sheet (prefix_name / prefix_url), placed right after
#    `mapping set reference`.
# ---------------------------------------------------------------------------
$mappingSetReference = $wb.Worksheets.Item("mapping set reference")
$prefix = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $mappingSetReference)
$prefix.Name = "prefix"
$prefix.Range("A1").Value = "prefix_name"
$prefix.Range("B1").Value = "prefix_url"

# ---------------------------------------------------------------------------
# 5. New, empty `NoTermFound` sheet, placed at the very end (after
#    `Propagatable`).
# ---------------------------------------------------------------------------
$propagatable = $wb.Worksheets.Item("Propagatable")
$noTermFound = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $propagatable)
$noTermFound.Name = "NoTermFound"
